# MODEL-INPUT CHANGES, removing lost prob transitions
#
# The "Transitions" sheet has a diagonal band of 1's (probability-of-loss
# markers) running through columns AH:AM on rows 16-33. These are being
# removed (cleared back to blank, formatting/style left intact) and the
# active sheet/selection moves from "Parameters" to "Transitions".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transitions")

# Cells whose <v>1</v> payload is removed (cell keeps its style, becomes
# an empty cell) - the "lost prob" diagonal entries in columns AH:AM.
$cellsToClear = @(
    "AH16", "AI16",
    "AI17", "AJ17",
    "AJ18", "AK18",
    "AK19", "AL19",
    "AL20", "AM20",
    "AM21",
    "AH22", "AI22",
    "AI23", "AJ23",
    "AI24", "AJ24", "AK24",
    "AJ25", "AK25", "AL25",
    "AK26", "AL26", "AM26",
    "AL27", "AM27",
    "AH28", "AI28",
    "AI29",
    "AI30", "AJ30",
    "AJ31", "AK31",
    "AK32", "AL32",
    "AL33", "AM33"
)

foreach ($ref in $cellsToClear) {
    $ws.Range($ref).Value = ""
}

# Move the active tab / selection from Parameters to Transitions.
$ws.Activate()
$ws.Range("P34").Select()
